# Update "want to go" counts (column F) for a few exhibition rows in both
# the "展览" and "全部类型" worksheets.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("F3").Value = 2539
    $ws.Range("F4").Value = 505
    $ws.Range("F6").Value = 6560
    $ws.Range("F7").Value = 388
}
